$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new LeetCode entries (rows 165-167) below the last existing row
# (164), mirroring the formatting of row 164 (style ids 1/2/2/1/1/1/1/4/4,
# wrapped text columns, centered alignment, date columns H/I).
# ---------------------------------------------------------------------------

function Add-ProblemRow($Row, $Num, $Name, $Tags, $Difficulty, $Success, $Fail, $Time, $FirstDate, $LastUpdate, $RowHeight) {

    # Copy the formatting of the template row (164) into the new row first so
    # that styles (borders/alignment/wrap/number formats) match exactly.
    $ws.Range("A164:I164").Copy() | Out-Null
    $target = "A" + $Row + ":I" + $Row
    $ws.Range($target).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($Row, 1).Value = $Num
    $ws.Cells.Item($Row, 2).Value = $Name
    $ws.Cells.Item($Row, 3).Value = $Tags
    $ws.Cells.Item($Row, 4).Value = $Difficulty
    $ws.Cells.Item($Row, 5).Value = $Success
    $ws.Cells.Item($Row, 6).Value = $Fail
    $ws.Cells.Item($Row, 7).Value = $Time
    $ws.Cells.Item($Row, 8).Value = $FirstDate
    $ws.Cells.Item($Row, 9).Value = $LastUpdate

    $ws.Rows.Item($Row).RowHeight = $RowHeight
}

Add-ProblemRow 165 2353 "Design a Food Rating System" `
    "#array #string #hash-table #set " "medium" `
    0 1 23 45917 45917 34

Add-ProblemRow 166 3005 "Count Elements With Maximum Frequency" `
    "#array #hash-table #counting " "easy" `
    1 0 5 45922 45922 34

Add-ProblemRow 167 1152 "Analyze User Website Visit Pattern" `
    "#array #hash-table #string #permutation #combination" "medium" `
    0 1 10 45922 45922 68

# Update the active selection to reflect where the user ended up editing.
$ws.Range("G171").Select() | Out-Null

Write-Host "Added rows 165-167"
